$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 205 — pushes the existing rows 205-220 down to 206-221
# (the weekly price-series gains one more week's record at the top of this
# block, Excel's default insert semantics shift everything below it down).
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with this week's record.
$ws.Range("A205").Value2 = 8
$ws.Range("B205").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C205").Value2 = "Coquimbo"
$ws.Range("D205").Value2 = 44769
$ws.Range("E205").Value2 = 4
$ws.Range("F205").Value2 = 100112037
$ws.Range("G205").Value2 = "Cebollín"
$ws.Range("H205").Value2 = "Sin especificar"
$ws.Range("I205").Value2 = "Primera"
$ws.Range("J205").Value2 = 1360
$ws.Range("K205").Value2 = 1400
$ws.Range("L205").Value2 = 1600
$ws.Range("M205").Value2 = 1500
$ws.Range("N205").Value2 = "`$/paquete 6 unidades"
$ws.Range("O205").Value2 = "Provincia del Elquí"
$ws.Range("P205").Value2 = 250
$ws.Range("Q205").Value2 = 6
$ws.Range("R205").Value2 = "Hortaliza"
